# Weekly price-sheet update: a new week of data is inserted as row 15,
# pushing the existing rows 15-46 down to 16-47 (dimension grows to A1:R47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15 (shifts old rows 15..46 -> 16..47,
# carrying their values/formatting down with them).
$ws.Rows.Item(15).Insert()

# The row that used to be 15 is now row 16; reuse its still-constant
# columns (market/region/category/quality/unit/origin/classification)
# for the brand-new row 15 instead of retyping the literals.
$ws.Range("A15").Value = $ws.Range("A16").Value()
$ws.Range("B15").Value = $ws.Range("B16").Value()
$ws.Range("C15").Value = $ws.Range("C16").Value()
$ws.Range("E15").Value = $ws.Range("E16").Value()
$ws.Range("F15").Value = $ws.Range("F16").Value()
$ws.Range("G15").Value = $ws.Range("G16").Value()
$ws.Range("H15").Value = $ws.Range("H16").Value()
$ws.Range("I15").Value = $ws.Range("I16").Value()
$ws.Range("N15").Value = $ws.Range("N16").Value()
$ws.Range("O15").Value = $ws.Range("O16").Value()
$ws.Range("Q15").Value = $ws.Range("Q16").Value()
$ws.Range("R15").Value = $ws.Range("R16").Value()

# New week's figures for the inserted row.
$ws.Range("D15").Value = "10/28/2021"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 750
$ws.Range("L15").Value = 800
$ws.Range("M15").Value = 775
$ws.Range("P15").Value = 258
